$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The due-date text "12, 2022" is retyped as four separate bold runs that
#    together read "11, 2023":
#       "1" | "1" | ", 202" | "3"
#    i.e. "Due: Monday, September 12, 2022" -> "Due: Monday, September 11, 2023"
# ---------------------------------------------------------------------------
$dateRng = $d.Content
$found = $dateRng.Find.Execute("12, 2022")
if (-not $found) {
    throw "Could not find date text '12, 2022' to replace."
}

$start = $dateRng.Start
$end = $dateRng.End

# Remove the original "12, 2022" text entirely, then retype it as separate runs.
$old = $d.Range($start, $end)
$old.Text = ""

$pieces = @("1", "1", ", 202", "3")
$pos = $start
foreach ($piece in $pieces) {
    $insertPoint = $d.Range($pos, $pos)
    $insertPoint.InsertAfter($piece)

    # Re-apply the run's bold formatting (collapsed ranges at a paragraph/run
    # boundary can otherwise pick up ambiguous/mixed formatting).
    $newRun = $d.Range($pos, $pos + $piece.Length)
    $newRun.Font.Bold = 1

    $pos = $pos + $piece.Length
}

# ---------------------------------------------------------------------------
# 2) Fix the typo "Show you're your work" -> "Show your work"
# ---------------------------------------------------------------------------
$apostrophe = [char]0x2019
$oldPhrase = "Show you" + $apostrophe + "re your work"
$newPhrase = "Show your work"

$typoRng = $d.Content
$typoFound = $typoRng.Find.Execute($oldPhrase, $true, $false, $false, $false, $false, $true, 1, $false, $newPhrase, 2)
if (-not $typoFound) {
    throw "Could not find typo text to replace."
}
